$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update header text (Volume number and report date range) ----
$ws.Range("A8").Value = "Volume 31   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/3/2024  Through  6/9/2024"

# ---- Helper anchors used to correctly convert cell types while preserving existing cell styles ----
# (I23 = numeric style "#,##0"; C29 = text placeholder "0"; E29 = text placeholder "***.*";
#  none of these three cells are themselves modified by this edit, so they are safe to use as format sources.)
$numAnchor = $ws.Range("I23")
$text0Anchor = $ws.Range("C29")
$textStarAnchor = $ws.Range("E29")

# ---- Convert cells that must become the text placeholder "0" (shared string, style 14) ----
$text0Anchor.Copy() | Out-Null
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("G22").PasteSpecial(-4163)
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("F33").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---- Convert cells that must become the text placeholder "***.*" (shared string, style 14) ----
$textStarAnchor.Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("H22").PasteSpecial(-4163)
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("E33").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# ---- Convert cells that are currently text placeholders but must become numeric (style 15: #,##0) ----
$numAnchor.Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C17").Value = 1
$ws.Range("C22").Value = 2

# ---- Update plain numeric cells (style already correct; only the value changes) ----
$ws.Range("F15").Value = 1
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = 6.060606060606
$ws.Range("L16").Value = 2.941176470588
$ws.Range("M16").Value = -20.454545454545
$ws.Range("N16").Value = -88.41059602649
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 41
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = -8.888888888888
$ws.Range("L17").Value = -6.818181818181
$ws.Range("M17").Value = 64
$ws.Range("N17").Value = -25.454545454545
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 40
$ws.Range("K18").Value = -27.272727272727
$ws.Range("L18").Value = -18.367346938775
$ws.Range("M18").Value = -25.925925925925
$ws.Range("N18").Value = -93.079584775086
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -53.333333333333
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 45
$ws.Range("H19").Value = -15.555555555555
$ws.Range("I19").Value = 180
$ws.Range("J19").Value = 216
$ws.Range("K19").Value = -16.666666666666
$ws.Range("L19").Value = -23.728813559322
$ws.Range("M19").Value = 18.421052631578
$ws.Range("N19").Value = -57.446808510638
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 25
$ws.Range("I20").Value = 75
$ws.Range("J20").Value = 53
$ws.Range("K20").Value = 41.509433962264
$ws.Range("L20").Value = 114.285714285714
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -95.052770448548
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -15.78947368421
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = 2.5
$ws.Range("I21").Value = 377
$ws.Range("J21").Value = 404
$ws.Range("K21").Value = -6.683168316831
$ws.Range("L21").Value = -7.142857142857
$ws.Range("M21").Value = 15.644171779141
$ws.Range("N21").Value = -86.900625434329
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = 63.636363636363
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = 125
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 45.833333333333
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 1.612903225806
$ws.Range("I24").Value = 700
$ws.Range("J24").Value = 698
$ws.Range("K24").Value = 0.286532951289
$ws.Range("L24").Value = -11.838790931989
$ws.Range("M24").Value = 70.316301703163
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -4.761904761904
$ws.Range("F25").Value = 84
$ws.Range("G25").Value = 90
$ws.Range("H25").Value = -6.666666666666
$ws.Range("I25").Value = 514
$ws.Range("J25").Value = 501
$ws.Range("K25").Value = 2.594810379241
$ws.Range("L25").Value = -11.072664359861
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 28
$ws.Range("I26").Value = 125
$ws.Range("J26").Value = 109
$ws.Range("K26").Value = 14.678899082568
$ws.Range("L26").Value = 52.439024390243
$ws.Range("M26").Value = 27.551020408163
$ws.Range("F27").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 14
$ws.Range("K28").Value = 7.692307692307
$ws.Range("L28").Value = -22.222222222222
$ws.Range("L31").Value = -37.5
$ws.Range("H33").Value = -100

